# Weekly fruit/vegetable data update: insert a new observation row right
# after row 294 (pushing the existing rows 295-317 down to 296-318) and
# populate it with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 295; this shifts rows 295:317 down to 296:318 and
# extends the used range to A1:R318.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new weekly record.
$ws.Range("A295").Value = 9
$ws.Range("B295").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C295").Value = "Metropolitana"
$ws.Range("D295").Value = 45041
$ws.Range("E295").Value = 13
$ws.Range("F295").Value = 100112026
$ws.Range("G295").Value = "Haba"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 52
$ws.Range("K295").Value = 20000
$ws.Range("L295").Value = 22000
$ws.Range("M295").Value = 21000
$ws.Range("N295").Value = "`$/saco 25 kilos"
$ws.Range("O295").Value = "Provincia de Limarí"
$ws.Range("P295").Value = 840
$ws.Range("Q295").Value = 25
$ws.Range("R295").Value = "Hortaliza"

Write-Output "row inserted and populated"
